$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "NR" scenario row (row 3) to the "TRMM" scenario
$ws.Range("D3").Value = "TRMM_population"
$ws.Range("F3").Value = "TRMM_population_count"

# Update the active selection to F3 (as reflected in the saved file)
$ws.Range("F3").Select()
